$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the preparer name typo "Retrofitted_0759" -> "S.GISH" in rows 23-25
$ws.Range("B23").Value = "S.GISH"
$ws.Range("B24").Value = "S.GISH"
$ws.Range("B25").Value = "S.GISH"

# Update the saved selection on the sheet to S21
$ws.Range("S21").Select()
